$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 14 (pushes old rows 14-18 down to 15-19,
# and the blank trailer rows 19-21 shift to 20-22).
$ws.Rows(14).Insert()

# Excel's default row-insert doesn't always carry the exact same
# direct-formatting as the row above for every column, so explicitly
# copy the formatting of row 13 onto the freshly inserted row 14.
$ws.Range("A13:P13").Copy()
$ws.Range("A14:P14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The insert pushed a formerly-blank placeholder row off the bottom
# (old row 21 -> new row 22); remove it so the used range stays A1:P21.
$ws.Rows(22).Delete()

# Re-sequence the "S.No." column for the new row and the rows that
# shifted down so it keeps counting 13, 14, 15 ... through row 19
# (which used to be row 18).
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18

# Populate the newly inserted row 14 with the weight / mean-square-error data.
$ws.Range("B14").Value = "April"
$ws.Range("C14").Value = 33
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 31.55
$ws.Range("F14").Value = 16.22
$ws.Range("G14").Value = 23.88
$ws.Range("H14").Value = 70
$ws.Range("I14").Value = 52.4
$ws.Range("J14").Value = 72.150000000000006
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = "Gurdaspur"
$ws.Range("M14").Value = "Sub Mountanious"
$ws.Range("N14").Value = "chari"
$ws.Range("O14").Value = "Wheat Straw"
$ws.Range("P14").Value = 16

# Update the selection to match the saved view (B2:P19, active cell B2).
$ws.Range("B2:P19").Select()
